$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: bump publish date/time
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true (must stay a literal text "true", not a Boolean,
# and must keep the same cell style as its neighbours)
$cell = $ws.Range("B17")
$cell.Value = "'true"
$ws.Range("B16").Copy()
$cell.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
